# Update cryptocurrency price/volume data per latest GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'76.200.49"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "'3.030.07"
$ws.Range("E3").Value = "  +3.40%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'197.20"
$ws.Range("E5").Value = "  -1.65%  "

$ws.Range("D6").Value = "'617.75"
$ws.Range("E6").Value = "  +3.54%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "'0.547"
$ws.Range("E8").Value = "  -1.08%  "

$ws.Range("E9").Value = "  +3.80%  "

$ws.Range("D10").Value = "'3.033.91"
$ws.Range("E10").Value = "  +3.53%  "

$ws.Range("D11").Value = "'0.437"
$ws.Range("E11").Value = "  -3.52%  "

$ws.Range("D12").Value = "'0.160"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").Value = "'5.24"
$ws.Range("E13").Value = "  +5.53%  "

$ws.Range("D14").Value = "'3.592.78"
$ws.Range("E14").Value = "  +3.45%  "

$ws.Range("D15").Value = "'28.77"
$ws.Range("E15").Value = "  +2.23%  "

$ws.Range("D16").Value = "'76.183.50"
$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").Value = "'0.0000192"
$ws.Range("E17").Value = "  +1.23%  "

$ws.Range("D18").Value = "'3.039.16"
$ws.Range("E18").Value = "  +3.24%  "

$ws.Range("D19").Value = "'13.48"
$ws.Range("E19").Value = "  +1.35%  "

$ws.Range("D20").Value = "'8.91"
$ws.Range("E20").Value = "  +1.71%  "

$ws.Range("D21").Value = "'380.08"
$ws.Range("E21").Value = "  +1.92%  "

$ws.Range("D22").Value = "'2.38"
$ws.Range("E22").Value = "  +3.52%  "

$ws.Range("D23").Value = "'4.35"
$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("D24").Value = "'3.176.37"
$ws.Range("E24").Value = "  +3.88%  "

$ws.Range("D25").Value = "'72.35"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").Value = "'4.32"
$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("D28").Value = "'9.75"
$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("D29").Value = "'0.0000107"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").Value = "'8.23"
$ws.Range("E31").Value = "  +1.83%  "

$ws.Range("D32").Value = "'1.38"
$ws.Range("E32").Value = "  -0.16%  "

$ws.Range("D33").Value = "'492.48"
$ws.Range("E33").Value = "  -0.71%  "

$ws.Range("E34").Value = "  +4.26%  "

$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("D36").Value = "'20.50"
$ws.Range("E36").Value = "  +1.48%  "

$ws.Range("D37").Value = "'162.24"
$ws.Range("E37").Value = "  -1.64%  "

$ws.Range("D38").Value = "'0.119"
$ws.Range("E38").Value = "  +8.20%  "

$ws.Range("D39").Value = "'20.02"
$ws.Range("E39").Value = "  +1.74%  "

$ws.Range("D40").Value = "'190.13"
$ws.Range("E40").Value = "  +6.27%  "

$ws.Range("D41").Value = "'0.378"
$ws.Range("E41").Value = "  -4.08%  "

$ws.Range("D42").Value = "'0.103"
$ws.Range("E42").Value = "  -6.14%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D44").Value = "'5.09"
$ws.Range("E44").Value = "  +3.00%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.778"
$ws.Range("E45").Value = "  +18.31%  "

$ws.Range("D46").Value = "'41.60"
$ws.Range("E46").Value = "  +3.54%  "

$ws.Range("D47").Value = "'1.24"
$ws.Range("E47").Value = "  +3.57%  "

$ws.Range("D48").Value = "'1.63"
$ws.Range("E48").Value = "  -1.73%  "

$ws.Range("D49").Value = "'2.42"
$ws.Range("E49").Value = "  +4.59%  "

$ws.Range("D50").Value = "'0.596"
$ws.Range("E50").Value = "  +1.20%  "

$ws.Range("D51").Value = "'3.85"
$ws.Range("E51").Value = "  -1.28%  "
